# Fix export to excel issue:
#  - insert two new header columns "主队排名" (after 主队) and "客队排名" (after 客队)
#  - append the three newly-scraped odds rows that were missing from the export
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
}

function Set-TextValues {
    param($ws, [hashtable]$values, [int]$row)
    foreach ($col in $values.Keys) {
        Set-TextValue $ws.Range("$col$row") $values[$col]
    }
    $cols = $values.Keys | Sort-Object { $ws.Range($_ + "1").Column }
    $first = $cols[0]
    $last = $cols[$cols.Count - 1]
    $rng = $ws.Range("$first" + "$row" + ":" + "$last" + "$row")
    $rng.Copy()
    $rng.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# 1. Insert "主队排名" right after "主队" (currently column G), and
#    "客队排名" right after "客队" (which, once the first column is
#    inserted, sits at column I).
# ---------------------------------------------------------------------------
$ws.Columns("G:G").Insert()
Set-TextValue $ws.Range("G1") "主队排名"
$ws.Range("G1").Copy()
$ws.Range("G1").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Columns("I:I").Insert()
Set-TextValue $ws.Range("I1") "客队排名"
$ws.Range("I1").Copy()
$ws.Range("I1").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Append the three scraped data rows.
# ---------------------------------------------------------------------------
$row2 = [ordered]@{
    A = "English Premier League"
    B = "370"
    C = "2017-2018"
    D = "2017-08-1"
    E = "2017,07-1,28,17,58,00"
    F = "Arsenal"
    G = "5"
    H = "Leicester City"
    I = "12"
    J = "阿森纳"
    K = "莱切斯特城"
    L = "70283107"
    M = "Oddset"
    N = "Oddset(德国)"
    O = "1.3"
    P = "4.25"
    Q = "6.5"
    U = "66.41"
    V = "20.31"
    W = "13.28"
}

$row3 = [ordered]@{
    A = "English Premier League"
    B = "82"
    C = "2017-2018"
    D = "2017-08-1"
    E = "2017,08-1,11,18,32,00"
    F = "Arsenal"
    G = "5"
    H = "Leicester City"
    I = "12"
    J = "阿森纳"
    K = "莱切斯特城"
    L = "69061424"
    M = "Ladbrokes"
    N = "立博(英国)"
    O = "1.36"
    P = "4.33"
    Q = "8"
    R = "1.4"
    S = "4.4"
    T = "6.5"
    U = "67.38"
    V = "21.16"
    W = "11.45"
    X = "65.21"
    Y = "20.75"
    Z = "14.04"
}

$row4 = [ordered]@{
    A = "English Premier League"
    B = "115"
    C = "2017-2018"
    D = "2017-08-1"
    E = "2017,08-1,11,18,40,00"
    F = "Arsenal"
    G = "5"
    H = "Leicester City"
    I = "12"
    J = "阿森纳"
    K = "莱切斯特城"
    L = "69060412"
    M = "William Hill"
    N = "威廉希尔(英国)"
    O = "1.53"
    P = "3.6"
    Q = "6"
    R = "1.44"
    S = "4.5"
    T = "7"
    U = "59.52"
    V = "25.3"
    W = "15.18"
    X = "65.54"
    Y = "20.97"
    Z = "13.48"
}

foreach ($col in $row2.Keys) { Set-TextValue $ws.Range("$col" + "2") $row2[$col] }
$ws.Range("A2:W2").Copy()
$ws.Range("A2:W2").PasteSpecial(-4163)
$excel.CutCopyMode = $false

foreach ($col in $row3.Keys) { Set-TextValue $ws.Range("$col" + "3") $row3[$col] }
$ws.Range("A3:Z3").Copy()
$ws.Range("A3:Z3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

foreach ($col in $row4.Keys) { Set-TextValue $ws.Range("$col" + "4") $row4[$col] }
$ws.Range("A4:Z4").Copy()
$ws.Range("A4:Z4").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Match the row height used by the header row on the new data rows.
# ---------------------------------------------------------------------------
$ws.Rows("2:4").RowHeight = 25
